# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
# A new weekly record is inserted at the top of the data block (row 179),
# pushing all existing records down by one row (179->180, ..., 288->289).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 179; everything below (old rows 179..288) shifts down
# to 180..289, carrying its data (and formatting) with it.
$ws.Rows("179:179").Insert()

# Populate the newly inserted row 179 with the new weekly record.
$ws.Cells.Item(179, 1).Value  = 10
$ws.Cells.Item(179, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(179, 3).Value  = "La Araucanía"
$ws.Cells.Item(179, 4).Value  = 45086
$ws.Cells.Item(179, 5).Value  = 9
$ws.Cells.Item(179, 6).Value  = "Fruta"
$ws.Cells.Item(179, 7).Value  = 100104
$ws.Cells.Item(179, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(179, 9).Value  = 100104003
$ws.Cells.Item(179, 10).Value = "Membrillo"
$ws.Cells.Item(179, 11).Value = "Champion"
$ws.Cells.Item(179, 12).Value = "Primera"
$ws.Cells.Item(179, 13).Value = 55
$ws.Cells.Item(179, 14).Value = 14000
$ws.Cells.Item(179, 15).Value = 14000
$ws.Cells.Item(179, 16).Value = 14000
$ws.Cells.Item(179, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(179, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(179, 19).Value = 778
$ws.Cells.Item(179, 20).Value = 18
